$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    ("Play Book of Relics for Free - Review and Ratings").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaStart = $metaPara.Range.Start
$insPoint = $d.Range($metaStart, $metaStart)
$insPoint.InsertAfter("Meta description: Read our unbiased review of Book of Relics, play for free and learn about its features, gameplay and betting range.")

# Bold only the "Meta description" label (the first 17 characters).
$labelLen = "Meta description".Length
$boldRange = $d.Range($metaStart, $metaStart + $labelLen)
$boldRange.Font.Bold = $true

# ------------------------------------------------------------------
# 2) Remove the duplicated bold heading paragraph near the end of the
#    document ("Play Book of Relics for Free - Review and Ratings").
# ------------------------------------------------------------------
$dupFound = $true
while ($dupFound) {
    $dupFound = $false
    for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text.TrimEnd() -eq "Play Book of Relics for Free - Review and Ratings") {
            $para.Range.Delete()
            $dupFound = $true
            break
        }
    }
}

# ------------------------------------------------------------------
# 3) Replace the trailing italic paragraph's text (formerly the meta
#    description, now duplicated at the top) with the AI image
#    generation prompt, keeping the paragraph's existing italic
#    formatting and straight (non-curly) quotation marks. Search
#    starting at that last paragraph specifically, since the same
#    sentence now also appears earlier (in the new meta-description
#    paragraph) and must be left untouched there.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$promptRange = $lastPara.Range
$promptRange.Find.Execute("Read our unbiased review of Book of Relics, play for free and learn about its features, gameplay and betting range.") | Out-Null
$promptRange.Text = 'Prompt: Create a cartoon-style feature image for "Book of Relics" that features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior standing in front of an ancient temple or pyramid, holding the Book of Relics in one hand and a handful of gold coins in the other. The background should be a desert landscape, with palm trees and sand dunes visible in the distance. The Maya warrior should be depicted with a big smile on their face, looking excited and confident as they hold their treasures. The image should evoke a sense of adventure, excitement, and the thrill of discovering ancient relics and treasure.'

Write-Host "done"
